$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.004940474577324
$ws.Range("C2").Value = 0.4825389026844107
$ws.Range("D2").Value = 0.02643071087631199
$ws.Range("F2").Value = 0.5153490970544112
$ws.Range("G2").Value = 0.3569480634963256
$ws.Range("H2").Value = 0.5256821673189833
$ws.Range("L2").Value = 0.2811251149203358
$ws.Range("M2").Value = 0.2330888419672235
$ws.Range("O2").Value = 1.703102221541783
$ws.Range("B3").Value = 0.8941744969332603
$ws.Range("C3").Value = 0.4720568239987415
$ws.Range("D3").Value = 0.02378328808529773
$ws.Range("F3").Value = 0.5164416939215926
$ws.Range("G3").Value = 0.359363130363711
$ws.Range("H3").Value = 0.5310352935070384
$ws.Range("L3").Value = 0.2802765279211528
$ws.Range("M3").Value = 0.2163400783165415
$ws.Range("O3").Value = 1.719153169090816
$ws.Range("B4").Value = 0.8260309281195077
$ws.Range("C4").Value = 0.4656974902366073
$ws.Range("D4").Value = 0.02214867525215425
$ws.Range("F4").Value = 0.517552616872706
$ws.Range("G4").Value = 0.3612184767894604
$ws.Range("H4").Value = 0.5346372054094033
$ws.Range("L4").Value = 0.2799190178454296
$ws.Range("M4").Value = 0.2060930472340203
$ws.Range("O4").Value = 1.73044817578824
$ws.Range("B5").Value = 0.7982303255030274
$ws.Range("C5").Value = 0.4631256482954313
$ws.Range("D5").Value = 0.02148031582805032
$ws.Range("F5").Value = 0.5181158488529221
$ws.Range("G5").Value = 0.3620679796257988
$ws.Range("H5").Value = 0.536184222067476
$ws.Range("L5").Value = 0.2798144830987823
$ws.Range("M5").Value = 0.2019268114732284
$ws.Range("O5").Value = 1.73541246270257
$ws.Range("B6").Value = 0.793612206135208
$ws.Range("C6").Value = 0.4626997923608087
$ws.Range("D6").Value = 0.02136920114585905
$ws.Range("F6").Value = 0.5182160443390842
$ws.Range("G6").Value = 0.3622146753341795
$ws.Range("H6").Value = 0.5364458864375905
$ws.Range("L6").Value = 0.2797996117472721
$ws.Range("M6").Value = 0.201235593561286
$ws.Range("O6").Value = 1.736258595542495
$ws.Range("B7").Value = 0.8256561247314949
$ws.Range("C7").Value = 0.4656627254965997
$ws.Range("D7").Value = 0.02213967054142074
$ws.Range("F7").Value = 0.5175597654985609
$ws.Range("G7").Value = 0.3612295554773937
$ws.Range("H7").Value = 0.5346577483312913
$ws.Range("L7").Value = 0.279917441378835
$ws.Range("M7").Value = 0.2060368210532673
$ws.Range("O7").Value = 1.730513662991228
$ws.Range("B8").Value = 0.9667770031406349
$ws.Range("C8").Value = 0.4789089670633473
$ws.Range("D8").Value = 0.02551978993834325
$ws.Range("F8").Value = 0.5156344063210483
$ws.Range("G8").Value = 0.3577033303084605
$ws.Range("H8").Value = 0.5274625210838906
$ws.Range("L8").Value = 0.2807986014307033
$ws.Range("M8").Value = 0.227306387956375
$ws.Range("O8").Value = 1.708337465146641
$ws.Range("B9").Value = 1.242393225368119
$ws.Range("C9").Value = 0.5054796365589311
$ws.Range("D9").Value = 0.03207453065946453
$ws.Range("F9").Value = 0.5153573558636708
$ws.Range("G9").Value = 0.3537546424750033
$ws.Range("H9").Value = 0.515853573020074
$ws.Range("L9").Value = 0.283823075510746
$ws.Range("M9").Value = 0.2692983298655847
$ws.Range("O9").Value = 1.676298025455779
$ws.Range("B10").Value = 1.444134200086921
$ws.Range("C10").Value = 0.5253475363590212
$ws.Range("D10").Value = 0.03684371631094052
$ws.Range("F10").Value = 0.5172971896808676
$ws.Range("G10").Value = 0.3526771481808879
$ws.Range("H10").Value = 0.5088504726383434
$ws.Range("L10").Value = 0.2868349655775688
$ws.Range("M10").Value = 0.300312526951366
$ws.Range("O10").Value = 1.659773462064067
$ws.Range("B11").Value = 1.535734062548954
$ws.Range("C11").Value = 0.5344579474733848
$ws.Range("D11").Value = 0.03900290146631846
$ws.Range("F11").Value = 0.5186473942796397
$ws.Range("G11").Value = 0.3525862403917301
$ws.Range("H11").Value = 0.5059962482691418
$ws.Range("L11").Value = 0.2883765762648238
$ws.Range("M11").Value = 0.3144552422845237
$ws.Range("O11").Value = 1.653786932576452
$ws.Range("B12").Value = 1.570394090440459
$ws.Range("C12").Value = 0.537917923425681
$ws.Range("D12").Value = 0.03981900429999996
$ws.Range("F12").Value = 0.5192261079309972
$ws.Range("G12").Value = 0.3526094814497469
$ws.Range("H12").Value = 0.5049631271392059
$ws.Range("L12").Value = 0.2889849826865003
$ws.Range("M12").Value = 0.3198154202498245
$ws.Range("O12").Value = 1.651740686057792
$ws.Range("B13").Value = 1.562930658555047
$ws.Range("C13").Value = 0.5371723137722881
$ws.Range("D13").Value = 0.03964331090632811
$ws.Range("F13").Value = 0.519098470581234
$ws.Range("G13").Value = 0.3526019073435265
$ws.Range("H13").Value = 0.5051835056465421
$ws.Range("L13").Value = 0.2888528564097754
$ws.Range("M13").Value = 0.3186608081487847
$ws.Range("O13").Value = 1.652171555250533
$ws.Range("B14").Value = 1.53858611322056
$ws.Range("C14").Value = 0.5347424021347251
$ws.Range("D14").Value = 0.03907007374215254
$ws.Range("F14").Value = 0.5186936533052631
$ws.Range("G14").Value = 0.3525869952893856
$ws.Range("H14").Value = 0.5059102963349176
$ws.Range("L14").Value = 0.2884261367260734
$ws.Range("M14").Value = 0.3148961359822309
$ws.Range("O14").Value = 1.653614159468873
$ws.Range("B15").Value = 1.52367082148487
$ws.Range("C15").Value = 0.5332553109918194
$ws.Range("D15").Value = 0.03871874830961985
$ws.Range("F15").Value = 0.5184544757750587
$ws.Range("G15").Value = 0.352585378530236
$ws.Range("H15").Value = 0.5063616912484719
$ws.Range("L15").Value = 0.2881679655769034
$ws.Range("M15").Value = 0.3125907613936718
$ws.Range("O15").Value = 1.654526559276576
$ws.Range("B16").Value = 1.438144282367603
$ws.Range("C16").Value = 0.524753576457698
$ws.Range("D16").Value = 0.03670239639686201
$ws.Range("F16").Value = 0.5172183766854488
$ws.Range("G16").Value = 0.3526911449297216
$ws.Range("H16").Value = 0.5090436773398181
$ws.Range("L16").Value = 0.2867376660837238
$ws.Range("M16").Value = 0.299388931531908
$ws.Range("O16").Value = 1.66019554871329
$ws.Range("B17").Value = 1.385630783286217
$ws.Range("C17").Value = 0.5195563313455693
$ws.Range("D17").Value = 0.0354627478884737
$ws.Range("F17").Value = 0.5165799877503972
$ws.Range("G17").Value = 0.3528584682511493
$ws.Range("H17").Value = 0.5107739275889145
$ws.Range("L17").Value = 0.2859041259087718
$ws.Range("M17").Value = 0.2912985969168815
$ws.Range("O17").Value = 1.664065755715427
$ws.Range("B18").Value = 1.355410228871051
$ws.Range("C18").Value = 0.5165738491167815
$ws.Range("D18").Value = 0.03474876425706697
$ws.Range("F18").Value = 0.516256822319157
$ws.Range("G18").Value = 0.3529922724237906
$ws.Range("H18").Value = 0.5118003237762991
$ws.Range("L18").Value = 0.2854408384610991
$ws.Range("M18").Value = 0.2866484862379153
$ws.Range("O18").Value = 1.666435824467371
$ws.Range("B19").Value = 1.345175349599117
$ws.Range("C19").Value = 0.5155652166644416
$ws.Range("D19").Value = 0.03450685639115392
$ws.Range("F19").Value = 0.5161549595553581
$ws.Range("G19").Value = 0.3530440201260419
$ws.Range("H19").Value = 0.5121532016598493
$ws.Range("L19").Value = 0.2852867506762209
$ws.Range("M19").Value = 0.2850746029607265
$ws.Range("O19").Value = 1.667263006781056
$ws.Range("B20").Value = 1.391222623571366
$ws.Range("C20").Value = 0.5201088820422228
$ws.Range("D20").Value = 0.03559481138036347
$ws.Range("F20").Value = 0.5166433883995936
$ws.Range("G20").Value = 0.3528367667925991
$ws.Range("H20").Value = 0.5105865099179852
$ws.Range("L20").Value = 0.2859911871629492
$ws.Range("M20").Value = 0.2921594939842151
$ws.Range("O20").Value = 1.663638854647814
$ws.Range("B21").Value = 1.545737440592632
$ws.Range("C21").Value = 0.53545585580639
$ws.Range("D21").Value = 0.03923848938886465
$ws.Range("F21").Value = 0.5188107269444231
$ws.Range("G21").Value = 0.3525898082712757
$ws.Range("H21").Value = 0.5056955250986306
$ws.Range("L21").Value = 0.2885508063595097
$ws.Range("M21").Value = 0.3160017877835557
$ws.Range("O21").Value = 1.653184436332481
$ws.Range("B22").Value = 1.64656439320089
$ws.Range("C22").Value = 0.5455444735282242
$ws.Range("D22").Value = 0.04161088236830324
$ws.Range("F22").Value = 0.520620257033535
$ws.Range("G22").Value = 0.3527646437627538
$ws.Range("H22").Value = 0.5027771020155001
$ws.Range("L22").Value = 0.2903672261447099
$ws.Range("M22").Value = 0.3316110138698249
$ws.Range("O22").Value = 1.6476386948992
$ws.Range("B23").Value = 1.592766197380513
$ws.Range("C23").Value = 0.5401547519941801
$ws.Range("D23").Value = 0.04034552690524151
$ws.Range("F23").Value = 0.5196184616575579
$ws.Range("G23").Value = 0.3526404816366409
$ws.Range("H23").Value = 0.5043092586063835
$ws.Range("L23").Value = 0.2893846420954702
$ws.Range("M23").Value = 0.3232777059073157
$ws.Range("O23").Value = 1.650480611156695
$ws.Range("B24").Value = 1.388694645393514
$ws.Range("C24").Value = 0.5198590567656538
$ws.Range("D24").Value = 0.03553510950063554
$ws.Range("F24").Value = 0.5166145883799871
$ws.Range("G24").Value = 0.3528464608962238
$ws.Range("H24").Value = 0.510671142854477
$ws.Range("L24").Value = 0.2859517771613582
$ws.Range("M24").Value = 0.2917702788781966
$ws.Range("O24").Value = 1.663831404892818
$ws.Range("B25").Value = 1.167959147567274
$ws.Range("C25").Value = 0.4982295753288213
$ws.Range("D25").Value = 0.03030935974905447
$ws.Range("F25").Value = 0.5150566593724548
$ws.Range("G25").Value = 0.3545037210874824
$ws.Range("H25").Value = 0.5187262200807652
$ws.Range("L25").Value = 0.2828661336895308
$ws.Range("M25").Value = 0.2579091476624455
$ws.Range("O25").Value = 1.683736234252279
